# Actualización automática 2025-05-30 16:20:08
# Apply updated sales figures to "VENTAS POR GRUPO" sheet of
# CASTRO ALCIVAR EDA MARIA workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# --- Row 5 ---
$ws.Range("K5").Value = 4132.31

# --- Row 6 ---
$ws.Range("D6").Value = 1139.71
$ws.Range("E6").Value = 69.42
$ws.Range("K6").Value = 313.21
$ws.Range("L6").Value = 6.05

# --- Row 11 ---
$ws.Range("M11").Value = 0

# --- Row 14 ---
$ws.Range("K14").Value = 353.28
$ws.Range("L14").Value = 684.29
$ws.Range("N14").Value = 1132.18

# --- Row 19 ---
$ws.Range("D19").Value = 0
$ws.Range("L19").Value = 0

# --- Row 22 ---
$ws.Range("C22").Value = 2052.86
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2088.56

# --- Row 24 ---
$ws.Range("D24").Value = 0
$ws.Range("L24").Value = 61.78

# --- Row 25 ---
$ws.Range("C25").Value = 0

# --- Row 26 ---
$ws.Range("C26").Value = 2140.99
$ws.Range("D26").Value = 864.96
$ws.Range("K26").Value = 1492.74
$ws.Range("L26").Value = 1751.07

# --- Row 27 ---
$ws.Range("D27").Value = 1205.95
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("L27").Value = 6109.34

# --- Row 28 ---
$ws.Range("D28").Value = 3218.72
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 678.51

# --- Row 29 ---
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = -26.21
$ws.Range("K29").Value = 2138.4
$ws.Range("L29").Value = 1451.1

# --- Row 31 ---
$ws.Range("K31").Value = 0

# --- Row 36 ---
$ws.Range("D36").Value = 518.98
$ws.Range("K36").Value = 2195.8
$ws.Range("L36").Value = 2727.15
$ws.Range("M36").Value = 262.99

# --- Row 39 ---
$ws.Range("D39").Value = 2403.41
$ws.Range("L39").Value = 0

# --- Row 41 ---
$ws.Range("L41").Value = 2689.09

# --- Row 42 ---
$ws.Range("D42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0

# --- Row 43 ---
$ws.Range("D43").Value = 0
$ws.Range("L43").Value = 4926.74
$ws.Range("M43").Value = 1327.27

# --- Row 44 ---
$ws.Range("D44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("L44").Value = 731.63

# --- Row 45 ---
$ws.Range("L45").Value = 722.54

# --- Row 47 ---
$ws.Range("G47").Value = 798
$ws.Range("H47").Value = 0

# --- Row 50 ---
$ws.Range("D50").Value = 1613.6
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 3339.53

# --- Row 51 ---
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = -11.75

# --- Row 54 ---
$ws.Range("C54").Value = 144

# --- Row 55 (summary counts "N de 53") ---
$ws.Range("C55").Value = "3 de 53"
$ws.Range("D55").Value = "7 de 53"
$ws.Range("G55").Value = "1 de 53"
$ws.Range("H55").Value = "0 de 53"
$ws.Range("K55").Value = "6 de 53"
$ws.Range("M55").Value = "2 de 53"
